$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory for check stock")

# Update the Location No. for row 2 from "01A012" to "01A011"
$ws.Range("F2").Value = "01A011"

# Update the remembered selection to match the author's last cursor position
$ws.Range("D12").Select()
